# Auto-generated Excel COM-interop script
# Applies the va_degree export fix values to the SCE short-circuit branch comparison workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("P2").Value = [double]"0"
$ws.Range("Q2").Value = [double]"0"
$ws.Range("P3").Value = [double]"0"
$ws.Range("Q3").Value = [double]"0"
$ws.Range("P4").Value = [double]"0"
$ws.Range("Q4").Value = [double]"0"
$ws.Range("P5").Value = [double]"0"
$ws.Range("Q5").Value = [double]"0"

$ws = $wb.Worksheets.Item(2)
$ws.Range("P2").Value = [double]"0"
$ws.Range("Q2").Value = [double]"0"
$ws.Range("P3").Value = [double]"0"
$ws.Range("Q3").Value = [double]"0"
$ws.Range("P4").Value = [double]"0"
$ws.Range("Q4").Value = [double]"0"
$ws.Range("P5").Value = [double]"0"
$ws.Range("Q5").Value = [double]"0"

$ws = $wb.Worksheets.Item(3)
$ws.Range("P2").Value = [double]"-10.05101633965328"
$ws.Range("Q2").Value = [double]"-9.529673900180338"
$ws.Range("P3").Value = [double]"-9.226445685756126"
$ws.Range("Q3").Value = [double]"-9.529673900180338"
$ws.Range("P4").Value = [double]"-9.226445685756126"
$ws.Range("Q4").Value = [double]"-9.529673900180338"
$ws.Range("P5").Value = [double]"-9.226445685756126"
$ws.Range("Q5").Value = [double]"-9.226445685756126"

$ws = $wb.Worksheets.Item(4)
$ws.Range("P2").Value = [double]"-10.05101633965328"
$ws.Range("Q2").Value = [double]"-9.529673900180338"
$ws.Range("P3").Value = [double]"-9.226445685756126"
$ws.Range("Q3").Value = [double]"-9.529673900180338"
$ws.Range("P4").Value = [double]"-9.226445685756126"
$ws.Range("Q4").Value = [double]"-9.529673900180338"
$ws.Range("P5").Value = [double]"-9.226445685756126"
$ws.Range("Q5").Value = [double]"-9.226445685756126"

$ws = $wb.Worksheets.Item(5)
$ws.Range("P2").Value = [double]"0"
$ws.Range("Q2").Value = [double]"0"
$ws.Range("P3").Value = [double]"0"
$ws.Range("Q3").Value = [double]"0"
$ws.Range("P4").Value = [double]"0"
$ws.Range("Q4").Value = [double]"0"
$ws.Range("P5").Value = [double]"0"
$ws.Range("Q5").Value = [double]"0"

$ws = $wb.Worksheets.Item(6)
$ws.Range("P2").Value = [double]"0"
$ws.Range("Q2").Value = [double]"0"
$ws.Range("P3").Value = [double]"0"
$ws.Range("Q3").Value = [double]"0"
$ws.Range("P4").Value = [double]"0"
$ws.Range("Q4").Value = [double]"0"
$ws.Range("P5").Value = [double]"0"
$ws.Range("Q5").Value = [double]"0"

$ws = $wb.Worksheets.Item(7)
$ws.Range("P2").Value = [double]"-10.4045966787725"
$ws.Range("Q2").Value = [double]"-7.634345070499246"
$ws.Range("P3").Value = [double]"-6.093293456603758"
$ws.Range("Q3").Value = [double]"-7.634345070499246"
$ws.Range("P4").Value = [double]"-6.093293456603762"
$ws.Range("Q4").Value = [double]"-7.634345070499246"
$ws.Range("P5").Value = [double]"-6.093293456603758"
$ws.Range("Q5").Value = [double]"-6.093293456603762"

$ws = $wb.Worksheets.Item(8)
$ws.Range("P2").Value = [double]"-10.4045966787725"
$ws.Range("Q2").Value = [double]"-7.634345070499246"
$ws.Range("P3").Value = [double]"-6.093293456603758"
$ws.Range("Q3").Value = [double]"-7.634345070499246"
$ws.Range("P4").Value = [double]"-6.093293456603762"
$ws.Range("Q4").Value = [double]"-7.634345070499246"
$ws.Range("P5").Value = [double]"-6.093293456603758"
$ws.Range("Q5").Value = [double]"-6.093293456603762"

$ws = $wb.Worksheets.Item(9)
$ws.Range("AL2").Value = [double]"1.591102290394586E-12"
$ws.Range("AM2").Value = [double]"179.9999999999773"
$ws.Range("AN2").Value = [double]"-179.9999999999874"
$ws.Range("AO2").Value = [double]"1.57040048766705E-12"
$ws.Range("AP2").Value = [double]"179.9999999999521"
$ws.Range("AQ2").Value = [double]"-179.9999999999613"
$ws.Range("AL3").Value = [double]"1.33885638337923E-12"
$ws.Range("AM3").Value = [double]"179.9999999999472"
$ws.Range("AN3").Value = [double]"-179.999999999955"
$ws.Range("AO3").Value = [double]"1.57040048766705E-12"
$ws.Range("AP3").Value = [double]"179.9999999999521"
$ws.Range("AQ3").Value = [double]"-179.9999999999613"
$ws.Range("AL4").Value = [double]"1.452591567033535E-12"
$ws.Range("AM4").Value = [double]"179.9999999999469"
$ws.Range("AN4").Value = [double]"-179.9999999999552"
$ws.Range("AO4").Value = [double]"1.57040048766705E-12"
$ws.Range("AP4").Value = [double]"179.9999999999521"
$ws.Range("AQ4").Value = [double]"-179.9999999999613"
$ws.Range("AL5").Value = [double]"1.33885638337923E-12"
$ws.Range("AM5").Value = [double]"179.9999999999472"
$ws.Range("AN5").Value = [double]"-179.999999999955"
$ws.Range("AO5").Value = [double]"1.452591567033535E-12"
$ws.Range("AP5").Value = [double]"179.9999999999469"
$ws.Range("AQ5").Value = [double]"-179.9999999999552"

$ws = $wb.Worksheets.Item(10)
$ws.Range("AL2").Value = [double]"1.591102290394586E-12"
$ws.Range("AM2").Value = [double]"179.9999999999773"
$ws.Range("AN2").Value = [double]"-179.9999999999874"
$ws.Range("AO2").Value = [double]"1.57040048766705E-12"
$ws.Range("AP2").Value = [double]"179.9999999999521"
$ws.Range("AQ2").Value = [double]"-179.9999999999613"
$ws.Range("AL3").Value = [double]"1.33885638337923E-12"
$ws.Range("AM3").Value = [double]"179.9999999999472"
$ws.Range("AN3").Value = [double]"-179.999999999955"
$ws.Range("AO3").Value = [double]"1.57040048766705E-12"
$ws.Range("AP3").Value = [double]"179.9999999999521"
$ws.Range("AQ3").Value = [double]"-179.9999999999613"
$ws.Range("AL4").Value = [double]"1.452591567033535E-12"
$ws.Range("AM4").Value = [double]"179.9999999999469"
$ws.Range("AN4").Value = [double]"-179.9999999999552"
$ws.Range("AO4").Value = [double]"1.57040048766705E-12"
$ws.Range("AP4").Value = [double]"179.9999999999521"
$ws.Range("AQ4").Value = [double]"-179.9999999999613"
$ws.Range("AL5").Value = [double]"1.33885638337923E-12"
$ws.Range("AM5").Value = [double]"179.9999999999472"
$ws.Range("AN5").Value = [double]"-179.999999999955"
$ws.Range("AO5").Value = [double]"1.452591567033535E-12"
$ws.Range("AP5").Value = [double]"179.9999999999469"
$ws.Range("AQ5").Value = [double]"-179.9999999999552"

$ws = $wb.Worksheets.Item(11)
$ws.Range("AL2").Value = [double]"6.057814235963239E-13"
$ws.Range("AM2").Value = [double]"-140.8543229646638"
$ws.Range("AN2").Value = [double]"122.4159013159209"
$ws.Range("AO2").Value = [double]"6.925037144414312E-13"
$ws.Range("AP2").Value = [double]"-142.867278420403"
$ws.Range("AQ2").Value = [double]"128.2389811279578"
$ws.Range("AL3").Value = [double]"6.538438094191336E-13"
$ws.Range("AM3").Value = [double]"-144.208653581376"
$ws.Range("AN3").Value = [double]"131.5664950534547"
$ws.Range("AO3").Value = [double]"6.925037144414312E-13"
$ws.Range("AP3").Value = [double]"-142.867278420403"
$ws.Range("AQ3").Value = [double]"128.2389811279578"
$ws.Range("AL4").Value = [double]"6.365009885364085E-13"
$ws.Range("AM4").Value = [double]"-144.208653581376"
$ws.Range("AN4").Value = [double]"131.5664950534547"
$ws.Range("AO4").Value = [double]"6.925037144414312E-13"
$ws.Range("AP4").Value = [double]"-142.867278420403"
$ws.Range("AQ4").Value = [double]"128.2389811279578"
$ws.Range("AL5").Value = [double]"6.538438094191336E-13"
$ws.Range("AM5").Value = [double]"-144.208653581376"
$ws.Range("AN5").Value = [double]"131.5664950534547"
$ws.Range("AO5").Value = [double]"6.365009885364085E-13"
$ws.Range("AP5").Value = [double]"-144.208653581376"
$ws.Range("AQ5").Value = [double]"131.5664950534547"

$ws = $wb.Worksheets.Item(12)
$ws.Range("AL2").Value = [double]"6.057814235963239E-13"
$ws.Range("AM2").Value = [double]"-140.8543229646638"
$ws.Range("AN2").Value = [double]"122.4159013159209"
$ws.Range("AO2").Value = [double]"6.925037144414312E-13"
$ws.Range("AP2").Value = [double]"-142.867278420403"
$ws.Range("AQ2").Value = [double]"128.2389811279578"
$ws.Range("AL3").Value = [double]"6.538438094191336E-13"
$ws.Range("AM3").Value = [double]"-144.208653581376"
$ws.Range("AN3").Value = [double]"131.5664950534547"
$ws.Range("AO3").Value = [double]"6.925037144414312E-13"
$ws.Range("AP3").Value = [double]"-142.867278420403"
$ws.Range("AQ3").Value = [double]"128.2389811279578"
$ws.Range("AL4").Value = [double]"6.365009885364085E-13"
$ws.Range("AM4").Value = [double]"-144.208653581376"
$ws.Range("AN4").Value = [double]"131.5664950534547"
$ws.Range("AO4").Value = [double]"6.925037144414312E-13"
$ws.Range("AP4").Value = [double]"-142.867278420403"
$ws.Range("AQ4").Value = [double]"128.2389811279578"
$ws.Range("AL5").Value = [double]"6.538438094191336E-13"
$ws.Range("AM5").Value = [double]"-144.208653581376"
$ws.Range("AN5").Value = [double]"131.5664950534547"
$ws.Range("AO5").Value = [double]"6.365009885364085E-13"
$ws.Range("AP5").Value = [double]"-144.208653581376"
$ws.Range("AQ5").Value = [double]"131.5664950534547"

$ws = $wb.Worksheets.Item(13)
$ws.Range("AL2").Value = [double]"7.917329727573314E-13"
$ws.Range("AM2").Value = [double]"-179.999999999902"
$ws.Range("AN2").Value = [double]"179.9999999998949"
$ws.Range("AO2").Value = [double]"5.369173410328852E-13"
$ws.Range("AP2").Value = [double]"-179.9999999994768"
$ws.Range("AQ2").Value = [double]"179.9999999994718"
$ws.Range("AL3").Value = [double]"3.682839488060667E-13"
$ws.Range("AM3").Value = [double]"-179.9999999994164"
$ws.Range("AN3").Value = [double]"179.9999999994127"
$ws.Range("AO3").Value = [double]"5.369173410328852E-13"
$ws.Range("AP3").Value = [double]"-179.9999999994768"
$ws.Range("AQ3").Value = [double]"179.9999999994718"
$ws.Range("AL4").Value = [double]"4.490569187661566E-13"
$ws.Range("AM4").Value = [double]"-179.9999999994165"
$ws.Range("AN4").Value = [double]"179.9999999994125"
$ws.Range("AO4").Value = [double]"5.369173410328852E-13"
$ws.Range("AP4").Value = [double]"-179.9999999994768"
$ws.Range("AQ4").Value = [double]"179.9999999994718"
$ws.Range("AL5").Value = [double]"3.682839488060667E-13"
$ws.Range("AM5").Value = [double]"-179.9999999994164"
$ws.Range("AN5").Value = [double]"179.9999999994127"
$ws.Range("AO5").Value = [double]"4.490569187661566E-13"
$ws.Range("AP5").Value = [double]"-179.9999999994165"
$ws.Range("AQ5").Value = [double]"179.9999999994125"

$ws = $wb.Worksheets.Item(14)
$ws.Range("AL2").Value = [double]"7.917329727573314E-13"
$ws.Range("AM2").Value = [double]"-179.999999999902"
$ws.Range("AN2").Value = [double]"179.9999999998949"
$ws.Range("AO2").Value = [double]"5.369173410328852E-13"
$ws.Range("AP2").Value = [double]"-179.9999999994768"
$ws.Range("AQ2").Value = [double]"179.9999999994718"
$ws.Range("AL3").Value = [double]"3.682839488060667E-13"
$ws.Range("AM3").Value = [double]"-179.9999999994164"
$ws.Range("AN3").Value = [double]"179.9999999994127"
$ws.Range("AO3").Value = [double]"5.369173410328852E-13"
$ws.Range("AP3").Value = [double]"-179.9999999994768"
$ws.Range("AQ3").Value = [double]"179.9999999994718"
$ws.Range("AL4").Value = [double]"4.490569187661566E-13"
$ws.Range("AM4").Value = [double]"-179.9999999994165"
$ws.Range("AN4").Value = [double]"179.9999999994125"
$ws.Range("AO4").Value = [double]"5.369173410328852E-13"
$ws.Range("AP4").Value = [double]"-179.9999999994768"
$ws.Range("AQ4").Value = [double]"179.9999999994718"
$ws.Range("AL5").Value = [double]"3.682839488060667E-13"
$ws.Range("AM5").Value = [double]"-179.9999999994164"
$ws.Range("AN5").Value = [double]"179.9999999994127"
$ws.Range("AO5").Value = [double]"4.490569187661566E-13"
$ws.Range("AP5").Value = [double]"-179.9999999994165"
$ws.Range("AQ5").Value = [double]"179.9999999994125"

$ws = $wb.Worksheets.Item(15)
$ws.Range("AL2").Value = [double]"3.726724087904993E-13"
$ws.Range("AM2").Value = [double]"-141.4823925534192"
$ws.Range("AN2").Value = [double]"122.8363393546082"
$ws.Range("AO2").Value = [double]"2.599745963568909E-13"
$ws.Range("AP2").Value = [double]"-142.5584264979141"
$ws.Range("AQ2").Value = [double]"131.69891753725"
$ws.Range("AL3").Value = [double]"2.119857453272499E-13"
$ws.Range("AM3").Value = [double]"-143.7604788262511"
$ws.Range("AN3").Value = [double]"136.2873245880946"
$ws.Range("AO3").Value = [double]"2.599745963568909E-13"
$ws.Range("AP3").Value = [double]"-142.5584264979141"
$ws.Range("AQ3").Value = [double]"131.69891753725"
$ws.Range("AL4").Value = [double]"2.396347970433406E-13"
$ws.Range("AM4").Value = [double]"-143.7604788262512"
$ws.Range("AN4").Value = [double]"136.2873245880946"
$ws.Range("AO4").Value = [double]"2.599745963568909E-13"
$ws.Range("AP4").Value = [double]"-142.5584264979141"
$ws.Range("AQ4").Value = [double]"131.69891753725"
$ws.Range("AL5").Value = [double]"2.119857453272499E-13"
$ws.Range("AM5").Value = [double]"-143.7604788262511"
$ws.Range("AN5").Value = [double]"136.2873245880946"
$ws.Range("AO5").Value = [double]"2.396347970433406E-13"
$ws.Range("AP5").Value = [double]"-143.7604788262512"
$ws.Range("AQ5").Value = [double]"136.2873245880946"

$ws = $wb.Worksheets.Item(16)
$ws.Range("AL2").Value = [double]"3.726724087904993E-13"
$ws.Range("AM2").Value = [double]"-141.4823925534192"
$ws.Range("AN2").Value = [double]"122.8363393546082"
$ws.Range("AO2").Value = [double]"2.599745963568909E-13"
$ws.Range("AP2").Value = [double]"-142.5584264979141"
$ws.Range("AQ2").Value = [double]"131.69891753725"
$ws.Range("AL3").Value = [double]"2.119857453272499E-13"
$ws.Range("AM3").Value = [double]"-143.7604788262511"
$ws.Range("AN3").Value = [double]"136.2873245880946"
$ws.Range("AO3").Value = [double]"2.599745963568909E-13"
$ws.Range("AP3").Value = [double]"-142.5584264979141"
$ws.Range("AQ3").Value = [double]"131.69891753725"
$ws.Range("AL4").Value = [double]"2.396347970433406E-13"
$ws.Range("AM4").Value = [double]"-143.7604788262512"
$ws.Range("AN4").Value = [double]"136.2873245880946"
$ws.Range("AO4").Value = [double]"2.599745963568909E-13"
$ws.Range("AP4").Value = [double]"-142.5584264979141"
$ws.Range("AQ4").Value = [double]"131.69891753725"
$ws.Range("AL5").Value = [double]"2.119857453272499E-13"
$ws.Range("AM5").Value = [double]"-143.7604788262511"
$ws.Range("AN5").Value = [double]"136.2873245880946"
$ws.Range("AO5").Value = [double]"2.396347970433406E-13"
$ws.Range("AP5").Value = [double]"-143.7604788262512"
$ws.Range("AQ5").Value = [double]"136.2873245880946"

$ws = $wb.Worksheets.Item(17)
$ws.Range("AL2").Value = [double]"0.05649667790076105"
$ws.Range("AN2").Value = [double]"0"
$ws.Range("AO2").Value = [double]"0.2021272668631934"
$ws.Range("AQ2").Value = [double]"0"
$ws.Range("AL3").Value = [double]"0.2765061667267382"
$ws.Range("AN3").Value = [double]"0"
$ws.Range("AO3").Value = [double]"0.2021272668631934"
$ws.Range("AQ3").Value = [double]"0"
$ws.Range("AL4").Value = [double]"0.2765061667267146"
$ws.Range("AN4").Value = [double]"0"
$ws.Range("AO4").Value = [double]"0.2021272668631934"
$ws.Range("AQ4").Value = [double]"0"
$ws.Range("AL5").Value = [double]"0.2765061667267382"
$ws.Range("AN5").Value = [double]"0"
$ws.Range("AO5").Value = [double]"0.2765061667267146"
$ws.Range("AQ5").Value = [double]"0"

$ws = $wb.Worksheets.Item(18)
$ws.Range("AL2").Value = [double]"0.05649667790076105"
$ws.Range("AN2").Value = [double]"0"
$ws.Range("AO2").Value = [double]"0.2021272668631934"
$ws.Range("AQ2").Value = [double]"0"
$ws.Range("AL3").Value = [double]"0.2765061667267382"
$ws.Range("AN3").Value = [double]"0"
$ws.Range("AO3").Value = [double]"0.2021272668631934"
$ws.Range("AQ3").Value = [double]"0"
$ws.Range("AL4").Value = [double]"0.2765061667267146"
$ws.Range("AN4").Value = [double]"0"
$ws.Range("AO4").Value = [double]"0.2021272668631934"
$ws.Range("AQ4").Value = [double]"0"
$ws.Range("AL5").Value = [double]"0.2765061667267382"
$ws.Range("AN5").Value = [double]"0"
$ws.Range("AO5").Value = [double]"0.2765061667267146"
$ws.Range("AQ5").Value = [double]"0"

$ws = $wb.Worksheets.Item(19)
$ws.Range("AL2").Value = [double]"2.05972647919493"
$ws.Range("AM2").Value = [double]"-106.7173761755005"
$ws.Range("AN2").Value = [double]"85.88603023776172"
$ws.Range("AO2").Value = [double]"1.922347956102813"
$ws.Range("AP2").Value = [double]"-106.6455958068003"
$ws.Range("AQ2").Value = [double]"86.86130100824833"
$ws.Range("AL3").Value = [double]"1.834400787235734"
$ws.Range("AM3").Value = [double]"-106.6064488821133"
$ws.Range("AN3").Value = [double]"87.43698758045056"
$ws.Range("AO3").Value = [double]"1.922347956102813"
$ws.Range("AP3").Value = [double]"-106.6455958068003"
$ws.Range("AQ3").Value = [double]"86.86130100824833"
$ws.Range("AL4").Value = [double]"1.834400787235711"
$ws.Range("AM4").Value = [double]"-106.6064488821133"
$ws.Range("AN4").Value = [double]"87.43698758045051"
$ws.Range("AO4").Value = [double]"1.922347956102813"
$ws.Range("AP4").Value = [double]"-106.6455958068003"
$ws.Range("AQ4").Value = [double]"86.86130100824833"
$ws.Range("AL5").Value = [double]"1.834400787235734"
$ws.Range("AM5").Value = [double]"-106.6064488821133"
$ws.Range("AN5").Value = [double]"87.43698758045056"
$ws.Range("AO5").Value = [double]"1.834400787235711"
$ws.Range("AP5").Value = [double]"-106.6064488821133"
$ws.Range("AQ5").Value = [double]"87.43698758045051"

$ws = $wb.Worksheets.Item(20)
$ws.Range("AL2").Value = [double]"2.05972647919493"
$ws.Range("AM2").Value = [double]"-106.7173761755005"
$ws.Range("AN2").Value = [double]"85.88603023776172"
$ws.Range("AO2").Value = [double]"1.922347956102813"
$ws.Range("AP2").Value = [double]"-106.6455958068003"
$ws.Range("AQ2").Value = [double]"86.86130100824833"
$ws.Range("AL3").Value = [double]"1.834400787235734"
$ws.Range("AM3").Value = [double]"-106.6064488821133"
$ws.Range("AN3").Value = [double]"87.43698758045056"
$ws.Range("AO3").Value = [double]"1.922347956102813"
$ws.Range("AP3").Value = [double]"-106.6455958068003"
$ws.Range("AQ3").Value = [double]"86.86130100824833"
$ws.Range("AL4").Value = [double]"1.834400787235711"
$ws.Range("AM4").Value = [double]"-106.6064488821133"
$ws.Range("AN4").Value = [double]"87.43698758045051"
$ws.Range("AO4").Value = [double]"1.922347956102813"
$ws.Range("AP4").Value = [double]"-106.6455958068003"
$ws.Range("AQ4").Value = [double]"86.86130100824833"
$ws.Range("AL5").Value = [double]"1.834400787235734"
$ws.Range("AM5").Value = [double]"-106.6064488821133"
$ws.Range("AN5").Value = [double]"87.43698758045056"
$ws.Range("AO5").Value = [double]"1.834400787235711"
$ws.Range("AP5").Value = [double]"-106.6064488821133"
$ws.Range("AQ5").Value = [double]"87.43698758045051"

$ws = $wb.Worksheets.Item(21)
$ws.Range("AL2").Value = [double]"0.05982627481619465"
$ws.Range("AN2").Value = [double]"0"
$ws.Range("AO2").Value = [double]"0.3273031178306479"
$ws.Range("AQ2").Value = [double]"0"
$ws.Range("AL3").Value = [double]"0.4531365616877374"
$ws.Range("AN3").Value = [double]"0"
$ws.Range("AO3").Value = [double]"0.3273031178306479"
$ws.Range("AQ3").Value = [double]"0"
$ws.Range("AL4").Value = [double]"0.4531365616877806"
$ws.Range("AN4").Value = [double]"0"
$ws.Range("AO4").Value = [double]"0.3273031178306479"
$ws.Range("AQ4").Value = [double]"0"
$ws.Range("AL5").Value = [double]"0.4531365616877374"
$ws.Range("AN5").Value = [double]"0"
$ws.Range("AO5").Value = [double]"0.4531365616877806"
$ws.Range("AQ5").Value = [double]"0"

$ws = $wb.Worksheets.Item(22)
$ws.Range("AL2").Value = [double]"0.05982627481619465"
$ws.Range("AN2").Value = [double]"0"
$ws.Range("AO2").Value = [double]"0.3273031178306479"
$ws.Range("AQ2").Value = [double]"0"
$ws.Range("AL3").Value = [double]"0.4531365616877374"
$ws.Range("AN3").Value = [double]"0"
$ws.Range("AO3").Value = [double]"0.3273031178306479"
$ws.Range("AQ3").Value = [double]"0"
$ws.Range("AL4").Value = [double]"0.4531365616877806"
$ws.Range("AN4").Value = [double]"0"
$ws.Range("AO4").Value = [double]"0.3273031178306479"
$ws.Range("AQ4").Value = [double]"0"
$ws.Range("AL5").Value = [double]"0.4531365616877374"
$ws.Range("AN5").Value = [double]"0"
$ws.Range("AO5").Value = [double]"0.4531365616877806"
$ws.Range("AQ5").Value = [double]"0"

$ws = $wb.Worksheets.Item(23)
$ws.Range("AL2").Value = [double]"2.055086913287214"
$ws.Range("AM2").Value = [double]"-107.1537816542741"
$ws.Range("AN2").Value = [double]"85.60853206981567"
$ws.Range("AO2").Value = [double]"1.84313717050446"
$ws.Range("AP2").Value = [double]"-104.9984074381874"
$ws.Range("AQ2").Value = [double]"89.0112927065375"
$ws.Range("AL3").Value = [double]"1.697274798967977"
$ws.Range("AM3").Value = [double]"-103.774646256571"
$ws.Range("AN3").Value = [double]"90.89841914161224"
$ws.Range("AO3").Value = [double]"1.84313717050446"
$ws.Range("AP3").Value = [double]"-104.9984074381874"
$ws.Range("AQ3").Value = [double]"89.0112927065375"
$ws.Range("AL4").Value = [double]"1.697274798967968"
$ws.Range("AM4").Value = [double]"-103.774646256571"
$ws.Range("AN4").Value = [double]"90.89841914161222"
$ws.Range("AO4").Value = [double]"1.84313717050446"
$ws.Range("AP4").Value = [double]"-104.9984074381874"
$ws.Range("AQ4").Value = [double]"89.0112927065375"
$ws.Range("AL5").Value = [double]"1.697274798967977"
$ws.Range("AM5").Value = [double]"-103.774646256571"
$ws.Range("AN5").Value = [double]"90.89841914161224"
$ws.Range("AO5").Value = [double]"1.697274798967968"
$ws.Range("AP5").Value = [double]"-103.774646256571"
$ws.Range("AQ5").Value = [double]"90.89841914161222"

$ws = $wb.Worksheets.Item(24)
$ws.Range("AL2").Value = [double]"2.055086913287214"
$ws.Range("AM2").Value = [double]"-107.1537816542741"
$ws.Range("AN2").Value = [double]"85.60853206981567"
$ws.Range("AO2").Value = [double]"1.84313717050446"
$ws.Range("AP2").Value = [double]"-104.9984074381874"
$ws.Range("AQ2").Value = [double]"89.0112927065375"
$ws.Range("AL3").Value = [double]"1.697274798967977"
$ws.Range("AM3").Value = [double]"-103.774646256571"
$ws.Range("AN3").Value = [double]"90.89841914161224"
$ws.Range("AO3").Value = [double]"1.84313717050446"
$ws.Range("AP3").Value = [double]"-104.9984074381874"
$ws.Range("AQ3").Value = [double]"89.0112927065375"
$ws.Range("AL4").Value = [double]"1.697274798967968"
$ws.Range("AM4").Value = [double]"-103.774646256571"
$ws.Range("AN4").Value = [double]"90.89841914161222"
$ws.Range("AO4").Value = [double]"1.84313717050446"
$ws.Range("AP4").Value = [double]"-104.9984074381874"
$ws.Range("AQ4").Value = [double]"89.0112927065375"
$ws.Range("AL5").Value = [double]"1.697274798967977"
$ws.Range("AM5").Value = [double]"-103.774646256571"
$ws.Range("AN5").Value = [double]"90.89841914161224"
$ws.Range("AO5").Value = [double]"1.697274798967968"
$ws.Range("AP5").Value = [double]"-103.774646256571"
$ws.Range("AQ5").Value = [double]"90.89841914161222"

$ws = $wb.Worksheets.Item(25)
$ws.Range("AL2").Value = [double]"0"
$ws.Range("AM2").Value = [double]"-147.7283134275611"
$ws.Range("AN2").Value = [double]"148.032989641989"
$ws.Range("AO2").Value = [double]"0"
$ws.Range("AP2").Value = [double]"-146.873188421504"
$ws.Range("AQ2").Value = [double]"147.9374049578102"
$ws.Range("AL3").Value = [double]"0"
$ws.Range("AM3").Value = [double]"-146.371367028733"
$ws.Range("AN3").Value = [double]"147.8040385391042"
$ws.Range("AO3").Value = [double]"0"
$ws.Range("AP3").Value = [double]"-146.873188421504"
$ws.Range("AQ3").Value = [double]"147.9374049578102"
$ws.Range("AL4").Value = [double]"0"
$ws.Range("AM4").Value = [double]"-146.371367028733"
$ws.Range("AN4").Value = [double]"147.8040385391042"
$ws.Range("AO4").Value = [double]"0"
$ws.Range("AP4").Value = [double]"-146.873188421504"
$ws.Range("AQ4").Value = [double]"147.9374049578102"
$ws.Range("AL5").Value = [double]"0"
$ws.Range("AM5").Value = [double]"-146.371367028733"
$ws.Range("AN5").Value = [double]"147.8040385391042"
$ws.Range("AO5").Value = [double]"0"
$ws.Range("AP5").Value = [double]"-146.371367028733"
$ws.Range("AQ5").Value = [double]"147.8040385391042"

$ws = $wb.Worksheets.Item(26)
$ws.Range("AL2").Value = [double]"0"
$ws.Range("AM2").Value = [double]"-147.7283134275611"
$ws.Range("AN2").Value = [double]"148.032989641989"
$ws.Range("AO2").Value = [double]"0"
$ws.Range("AP2").Value = [double]"-146.873188421504"
$ws.Range("AQ2").Value = [double]"147.9374049578102"
$ws.Range("AL3").Value = [double]"0"
$ws.Range("AM3").Value = [double]"-146.371367028733"
$ws.Range("AN3").Value = [double]"147.8040385391042"
$ws.Range("AO3").Value = [double]"0"
$ws.Range("AP3").Value = [double]"-146.873188421504"
$ws.Range("AQ3").Value = [double]"147.9374049578102"
$ws.Range("AL4").Value = [double]"0"
$ws.Range("AM4").Value = [double]"-146.371367028733"
$ws.Range("AN4").Value = [double]"147.8040385391042"
$ws.Range("AO4").Value = [double]"0"
$ws.Range("AP4").Value = [double]"-146.873188421504"
$ws.Range("AQ4").Value = [double]"147.9374049578102"
$ws.Range("AL5").Value = [double]"0"
$ws.Range("AM5").Value = [double]"-146.371367028733"
$ws.Range("AN5").Value = [double]"147.8040385391042"
$ws.Range("AO5").Value = [double]"0"
$ws.Range("AP5").Value = [double]"-146.371367028733"
$ws.Range("AQ5").Value = [double]"147.8040385391042"

$ws = $wb.Worksheets.Item(27)
$ws.Range("AL2").Value = [double]"-33.12503482398865"
$ws.Range("AM2").Value = [double]"-139.212415080177"
$ws.Range("AN2").Value = [double]"148.1841785064804"
$ws.Range("AO2").Value = [double]"-29.72036582402619"
$ws.Range("AP2").Value = [double]"-139.1548791201984"
$ws.Range("AQ2").Value = [double]"147.3382240646751"
$ws.Range("AL3").Value = [double]"-27.77641377061646"
$ws.Range("AM3").Value = [double]"-139.1268703916235"
$ws.Range("AN3").Value = [double]"146.8296931457364"
$ws.Range("AO3").Value = [double]"-29.72036582402619"
$ws.Range("AP3").Value = [double]"-139.1548791201984"
$ws.Range("AQ3").Value = [double]"147.3382240646751"
$ws.Range("AL4").Value = [double]"-27.77641377061655"
$ws.Range("AM4").Value = [double]"-139.1268703916235"
$ws.Range("AN4").Value = [double]"146.8296931457365"
$ws.Range("AO4").Value = [double]"-29.72036582402619"
$ws.Range("AP4").Value = [double]"-139.1548791201984"
$ws.Range("AQ4").Value = [double]"147.3382240646751"
$ws.Range("AL5").Value = [double]"-27.77641377061646"
$ws.Range("AM5").Value = [double]"-139.1268703916235"
$ws.Range("AN5").Value = [double]"146.8296931457364"
$ws.Range("AO5").Value = [double]"-27.77641377061655"
$ws.Range("AP5").Value = [double]"-139.1268703916235"
$ws.Range("AQ5").Value = [double]"146.8296931457365"

$ws = $wb.Worksheets.Item(28)
$ws.Range("AL2").Value = [double]"-33.12503482398865"
$ws.Range("AM2").Value = [double]"-139.212415080177"
$ws.Range("AN2").Value = [double]"148.1841785064804"
$ws.Range("AO2").Value = [double]"-29.72036582402619"
$ws.Range("AP2").Value = [double]"-139.1548791201984"
$ws.Range("AQ2").Value = [double]"147.3382240646751"
$ws.Range("AL3").Value = [double]"-27.77641377061646"
$ws.Range("AM3").Value = [double]"-139.1268703916235"
$ws.Range("AN3").Value = [double]"146.8296931457364"
$ws.Range("AO3").Value = [double]"-29.72036582402619"
$ws.Range("AP3").Value = [double]"-139.1548791201984"
$ws.Range("AQ3").Value = [double]"147.3382240646751"
$ws.Range("AL4").Value = [double]"-27.77641377061655"
$ws.Range("AM4").Value = [double]"-139.1268703916235"
$ws.Range("AN4").Value = [double]"146.8296931457365"
$ws.Range("AO4").Value = [double]"-29.72036582402619"
$ws.Range("AP4").Value = [double]"-139.1548791201984"
$ws.Range("AQ4").Value = [double]"147.3382240646751"
$ws.Range("AL5").Value = [double]"-27.77641377061646"
$ws.Range("AM5").Value = [double]"-139.1268703916235"
$ws.Range("AN5").Value = [double]"146.8296931457364"
$ws.Range("AO5").Value = [double]"-27.77641377061655"
$ws.Range("AP5").Value = [double]"-139.1268703916235"
$ws.Range("AQ5").Value = [double]"146.8296931457365"

$ws = $wb.Worksheets.Item(29)
$ws.Range("AL2").Value = [double]"0"
$ws.Range("AM2").Value = [double]"-147.6238196864016"
$ws.Range("AN2").Value = [double]"147.9448907708998"
$ws.Range("AO2").Value = [double]"0"
$ws.Range("AP2").Value = [double]"-146.3957340254715"
$ws.Range("AQ2").Value = [double]"148.1055649584091"
$ws.Range("AL3").Value = [double]"0"
$ws.Range("AM3").Value = [double]"-145.7020323973445"
$ws.Range("AN3").Value = [double]"148.0232122979875"
$ws.Range("AO3").Value = [double]"0"
$ws.Range("AP3").Value = [double]"-146.3957340254715"
$ws.Range("AQ3").Value = [double]"148.1055649584091"
$ws.Range("AL4").Value = [double]"0"
$ws.Range("AM4").Value = [double]"-145.7020323973445"
$ws.Range("AN4").Value = [double]"148.0232122979875"
$ws.Range("AO4").Value = [double]"0"
$ws.Range("AP4").Value = [double]"-146.3957340254715"
$ws.Range("AQ4").Value = [double]"148.1055649584091"
$ws.Range("AL5").Value = [double]"0"
$ws.Range("AM5").Value = [double]"-145.7020323973445"
$ws.Range("AN5").Value = [double]"148.0232122979875"
$ws.Range("AO5").Value = [double]"0"
$ws.Range("AP5").Value = [double]"-145.7020323973445"
$ws.Range("AQ5").Value = [double]"148.0232122979875"

$ws = $wb.Worksheets.Item(30)
$ws.Range("AL2").Value = [double]"0"
$ws.Range("AM2").Value = [double]"-147.6238196864016"
$ws.Range("AN2").Value = [double]"147.9448907708998"
$ws.Range("AO2").Value = [double]"0"
$ws.Range("AP2").Value = [double]"-146.3957340254715"
$ws.Range("AQ2").Value = [double]"148.1055649584091"
$ws.Range("AL3").Value = [double]"0"
$ws.Range("AM3").Value = [double]"-145.7020323973445"
$ws.Range("AN3").Value = [double]"148.0232122979875"
$ws.Range("AO3").Value = [double]"0"
$ws.Range("AP3").Value = [double]"-146.3957340254715"
$ws.Range("AQ3").Value = [double]"148.1055649584091"
$ws.Range("AL4").Value = [double]"0"
$ws.Range("AM4").Value = [double]"-145.7020323973445"
$ws.Range("AN4").Value = [double]"148.0232122979875"
$ws.Range("AO4").Value = [double]"0"
$ws.Range("AP4").Value = [double]"-146.3957340254715"
$ws.Range("AQ4").Value = [double]"148.1055649584091"
$ws.Range("AL5").Value = [double]"0"
$ws.Range("AM5").Value = [double]"-145.7020323973445"
$ws.Range("AN5").Value = [double]"148.0232122979875"
$ws.Range("AO5").Value = [double]"0"
$ws.Range("AP5").Value = [double]"-145.7020323973445"
$ws.Range("AQ5").Value = [double]"148.0232122979875"

$ws = $wb.Worksheets.Item(31)
$ws.Range("AL2").Value = [double]"-33.13110206389089"
$ws.Range("AM2").Value = [double]"-139.1611711215787"
$ws.Range("AN2").Value = [double]"148.0765151544961"
$ws.Range("AO2").Value = [double]"-26.31331180194492"
$ws.Range("AP2").Value = [double]"-139.1220723182847"
$ws.Range("AQ2").Value = [double]"146.8682275371679"
$ws.Range("AL3").Value = [double]"-22.54037501752185"
$ws.Range("AM3").Value = [double]"-139.1568268946852"
$ws.Range("AN3").Value = [double]"146.1701077073158"
$ws.Range("AO3").Value = [double]"-26.31331180194492"
$ws.Range("AP3").Value = [double]"-139.1220723182847"
$ws.Range("AQ3").Value = [double]"146.8682275371679"
$ws.Range("AL4").Value = [double]"-22.54037501752187"
$ws.Range("AM4").Value = [double]"-139.1568268946852"
$ws.Range("AN4").Value = [double]"146.1701077073158"
$ws.Range("AO4").Value = [double]"-26.31331180194492"
$ws.Range("AP4").Value = [double]"-139.1220723182847"
$ws.Range("AQ4").Value = [double]"146.8682275371679"
$ws.Range("AL5").Value = [double]"-22.54037501752185"
$ws.Range("AM5").Value = [double]"-139.1568268946852"
$ws.Range("AN5").Value = [double]"146.1701077073158"
$ws.Range("AO5").Value = [double]"-22.54037501752187"
$ws.Range("AP5").Value = [double]"-139.1568268946852"
$ws.Range("AQ5").Value = [double]"146.1701077073158"

$ws = $wb.Worksheets.Item(32)
$ws.Range("AL2").Value = [double]"-33.13110206389089"
$ws.Range("AM2").Value = [double]"-139.1611711215787"
$ws.Range("AN2").Value = [double]"148.0765151544961"
$ws.Range("AO2").Value = [double]"-26.31331180194492"
$ws.Range("AP2").Value = [double]"-139.1220723182847"
$ws.Range("AQ2").Value = [double]"146.8682275371679"
$ws.Range("AL3").Value = [double]"-22.54037501752185"
$ws.Range("AM3").Value = [double]"-139.1568268946852"
$ws.Range("AN3").Value = [double]"146.1701077073158"
$ws.Range("AO3").Value = [double]"-26.31331180194492"
$ws.Range("AP3").Value = [double]"-139.1220723182847"
$ws.Range("AQ3").Value = [double]"146.8682275371679"
$ws.Range("AL4").Value = [double]"-22.54037501752187"
$ws.Range("AM4").Value = [double]"-139.1568268946852"
$ws.Range("AN4").Value = [double]"146.1701077073158"
$ws.Range("AO4").Value = [double]"-26.31331180194492"
$ws.Range("AP4").Value = [double]"-139.1220723182847"
$ws.Range("AQ4").Value = [double]"146.8682275371679"
$ws.Range("AL5").Value = [double]"-22.54037501752185"
$ws.Range("AM5").Value = [double]"-139.1568268946852"
$ws.Range("AN5").Value = [double]"146.1701077073158"
$ws.Range("AO5").Value = [double]"-22.54037501752187"
$ws.Range("AP5").Value = [double]"-139.1568268946852"
$ws.Range("AQ5").Value = [double]"146.1701077073158"
